# Ausstände CA.xlsx - "Verlust Vogesen und Trikot Set Geld Hape eingetragen"
#
# The "Lutz" sheet had its expense rows (2-6) reordered (reversed) and two
# new expense rows were appended: "Verlust Vogesen" (-284.15) and
# "Trikot Set Hape Bar eingesackt" (-100). The currency number format used
# in column C was also changed from the accounting style to a simpler
# format that shows negative values in red.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lutz")

$newCurrencyFormat = "#,##0.00\ [$€-407];[Red]\-#,##0.00\ [$€-407]"
$mmmYyFormat = "mmm-yy"    # builtin numFmtId 14 -> style already used for 42736/42614
$shortDateFormat = "m/d/yy" # builtin numFmtId 17 -> style already used for 42707

# --- Re-enter rows 2-6 in their new (reversed) order -----------------------
$ws.Range("A2").Value = 42614
$ws.Range("A2").NumberFormat = $mmmYyFormat
$ws.Range("B2").Value = "Marketing, Recht, Gründung, Kauf und Verkauf eines Reisebüros"
$ws.Range("C2").Value = 138

$ws.Range("A3").Value = 42707
$ws.Range("A3").NumberFormat = $shortDateFormat
$ws.Range("B3").Value = "Plugin Instagram Feed Developer"
$ws.Range("C3").Value = 60

$ws.Range("A4").Value = 42736
$ws.Range("A4").NumberFormat = $mmmYyFormat
$ws.Range("B4").Value = "Plugin Hover Effects Builder"
$ws.Range("C4").Value = 11

$ws.Range("A5").ClearFormats()
$ws.Range("A5").Value = "März/April 2017"
$ws.Range("B5").Value = "Verpflegung Ligurien 7 Tage"
$ws.Range("C5").Value = 140

$ws.Range("A6").ClearFormats()
$ws.Range("A6").Value = "März/April 2017"
$ws.Range("B6").Value = "Guiding 3 Wochen"
$ws.Range("C6").Value = 735

# --- Append the two new expense rows ---------------------------------------
$ws.Range("A7").Value = 42856
$ws.Range("A7").NumberFormat = $mmmYyFormat
$ws.Range("B7").Value = "Verlust Vogesen"
$ws.Range("C7").Value = -284.14999999999998

$ws.Range("A8").Value = 42856
$ws.Range("A8").NumberFormat = $mmmYyFormat
$ws.Range("B8").Value = "Trikot Set Hape Bar eingesackt"
$ws.Range("C8").Value = -100

# --- Apply the updated currency number format to the whole amount column ---
$ws.Range("C1:C8").NumberFormat = $newCurrencyFormat

# --- Update the selected cell shown when the sheet is reopened -------------
$ws.Activate()
$ws.Range("F12").Select()

$wb.Save()
